$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value2 = 942.5
$ws.Cells.Item(18, 9).Value2 = 1093
$ws.Cells.Item(18, 11).Value2 = 1093
$ws.Cells.Item(18, 13).Value2 = -809
$ws.Cells.Item(19, 8).Value2 = 1333.625
$ws.Cells.Item(19, 9).Value2 = 1496.8334
$ws.Cells.Item(19, 11).Value2 = 1496.8334
$ws.Cells.Item(19, 13).Value2 = -1321.8334
$ws.Cells.Item(62, 8).Value2 = 38098750
$ws.Cells.Item(62, 9).Value2 = 44447730
$ws.Cells.Item(62, 10).Value2 = 4903
$ws.Cells.Item(62, 11).Value2 = 44447730
$ws.Cells.Item(62, 12).Value2 = 4903
$ws.Cells.Item(62, 13).Value2 = -44447106
$ws.Cells.Item(62, 14).Value2 = -6151
$ws.Cells.Item(65, 8).Value2 = 38098750
$ws.Cells.Item(65, 9).Value2 = 44447730
$ws.Cells.Item(65, 10).Value2 = 4903
$ws.Cells.Item(65, 11).Value2 = 222238650
$ws.Cells.Item(65, 12).Value2 = 24515
$ws.Cells.Item(65, 13).Value2 = -222235530
$ws.Cells.Item(65, 14).Value2 = -30755
$ws.Cells.Item(98, 8).Value2 = 1564.2963
$ws.Cells.Item(98, 9).Value2 = 1564.2963
$ws.Cells.Item(98, 11).Value2 = 1564.2963
$ws.Cells.Item(98, 13).Value2 = -66.29629999999997
$ws.Cells.Item(122, 8).Value2 = 1564.2963
$ws.Cells.Item(122, 9).Value2 = 1564.2963
$ws.Cells.Item(122, 11).Value2 = 4692.8889
$ws.Cells.Item(122, 13).Value2 = -2242.8889
$ws.Cells.Item(132, 8).Value2 = 1083.02
$ws.Cells.Item(132, 9).Value2 = 899.05
$ws.Cells.Item(132, 11).Value2 = 2697.15
$ws.Cells.Item(132, 13).Value2 = -167.1499999999996
$ws.Cells.Item(137, 8).Value2 = 16676287
$ws.Cells.Item(137, 10).Value2 = 19525
$ws.Cells.Item(137, 12).Value2 = 58575
$ws.Cells.Item(137, 14).Value2 = -63675

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value2 = 49360.22
$ws.Cells.Item(32, 9).Value2 = 53775.43
$ws.Cells.Item(32, 11).Value2 = 53775.43
$ws.Cells.Item(32, 13).Value2 = -53488.43
$ws.Cells.Item(132, 8).Value2 = 6430.3335
$ws.Cells.Item(132, 9).Value2 = 4787.9165
$ws.Cells.Item(132, 10).Value2 = 13000
$ws.Cells.Item(132, 11).Value2 = 14363.7495
$ws.Cells.Item(132, 12).Value2 = 39000
$ws.Cells.Item(132, 13).Value2 = -11833.7495
$ws.Cells.Item(132, 14).Value2 = -44060

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(97, 8).Value2 = 14699.4
$ws.Cells.Item(97, 9).Value2 = 10613.667
$ws.Cells.Item(97, 11).Value2 = 10613.667
$ws.Cells.Item(97, 13).Value2 = -9622.666999999999
$ws.Cells.Item(102, 8).Value2 = 13539.333
$ws.Cells.Item(102, 9).Value2 = 5233.143
$ws.Cells.Item(102, 11).Value2 = 5233.143
$ws.Cells.Item(102, 13).Value2 = -1988.143
$ws.Cells.Item(132, 8).Value2 = 98836
$ws.Cells.Item(132, 10).Value2 = 98836
$ws.Cells.Item(132, 12).Value2 = 98836
$ws.Cells.Item(132, 14).Value2 = -108956

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value2 = 34486620
$ws.Cells.Item(31, 9).Value2 = 142858860
$ws.Cells.Item(31, 11).Value2 = 142858860
$ws.Cells.Item(31, 13).Value2 = -142858565
$ws.Cells.Item(34, 8).Value2 = 34486620
$ws.Cells.Item(34, 9).Value2 = 142858860
$ws.Cells.Item(34, 11).Value2 = 142858860
$ws.Cells.Item(34, 13).Value2 = -142858658
$ws.Cells.Item(48, 8).Value2 = 39974.5
$ws.Cells.Item(48, 9).Value2 = 0
$ws.Cells.Item(48, 10).Value2 = 39974.5
$ws.Cells.Item(48, 11).Value2 = 0
$ws.Cells.Item(48, 13).Value2 = 39974.5
$ws.Cells.Item(48, 14).Value2 = -40926.5
$ws.Cells.Item(58, 8).Value2 = 4378.905
$ws.Cells.Item(58, 9).Value2 = 3160.6365
$ws.Cells.Item(58, 10).Value2 = 5719
$ws.Cells.Item(58, 11).Value2 = 3160.6365
$ws.Cells.Item(58, 12).Value2 = 5719
$ws.Cells.Item(58, 13).Value2 = -2957.6365
$ws.Cells.Item(58, 14).Value2 = -6125
$ws.Cells.Item(105, 8).Value2 = 2090
$ws.Cells.Item(105, 9).Value2 = 2160
$ws.Cells.Item(105, 10).Value2 = 1950
$ws.Cells.Item(105, 11).Value2 = 2160
$ws.Cells.Item(105, 12).Value2 = 1950
$ws.Cells.Item(105, 13).Value2 = -413
$ws.Cells.Item(105, 14).Value2 = -5444
$ws.Cells.Item(107, 8).Value2 = 1817.8096
$ws.Cells.Item(107, 10).Value2 = 2066.5833
$ws.Cells.Item(107, 12).Value2 = 2066.5833
$ws.Cells.Item(107, 14).Value2 = -5906.5833
$ws.Cells.Item(136, 8).Value2 = 4378.905
$ws.Cells.Item(136, 9).Value2 = 3160.6365
$ws.Cells.Item(136, 10).Value2 = 5719
$ws.Cells.Item(136, 11).Value2 = 9481.9095
$ws.Cells.Item(136, 12).Value2 = 17157
$ws.Cells.Item(136, 13).Value2 = -6931.9095
$ws.Cells.Item(136, 14).Value2 = -22257
$ws.Cells.Item(48, 12).ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value2 = 3290.5
$ws.Cells.Item(2, 10).Value2 = 5949
$ws.Cells.Item(2, 12).Value2 = 35694
$ws.Cells.Item(2, 14).Value2 = -35920
$ws.Cells.Item(38, 8).Value2 = 83
$ws.Cells.Item(38, 9).Value2 = 83
$ws.Cells.Item(38, 11).Value2 = 249
$ws.Cells.Item(38, 13).Value2 = 98
$ws.Cells.Item(75, 8).Value2 = 1682.6
$ws.Cells.Item(75, 9).Value2 = 1253
$ws.Cells.Item(75, 10).Value2 = 1969
$ws.Cells.Item(75, 11).Value2 = 3759
$ws.Cells.Item(75, 12).Value2 = 5907
$ws.Cells.Item(75, 13).Value2 = -2761
$ws.Cells.Item(75, 14).Value2 = -7903
$ws.Cells.Item(78, 8).Value2 = 1682.6
$ws.Cells.Item(78, 9).Value2 = 1253
$ws.Cells.Item(78, 10).Value2 = 1969
$ws.Cells.Item(78, 11).Value2 = 11277
$ws.Cells.Item(78, 12).Value2 = 17721
$ws.Cells.Item(78, 13).Value2 = -6285
$ws.Cells.Item(78, 14).Value2 = -27705
$ws.Cells.Item(88, 8).Value2 = 3950
$ws.Cells.Item(88, 9).Value2 = 3950
$ws.Cells.Item(88, 10).Value2 = 0
$ws.Cells.Item(88, 11).Value2 = 11850
$ws.Cells.Item(88, 12).Value2 = 0
$ws.Cells.Item(88, 14).Value2 = -11422
$ws.Cells.Item(91, 8).Value2 = 3950
$ws.Cells.Item(91, 9).Value2 = 3950
$ws.Cells.Item(91, 10).Value2 = 0
$ws.Cells.Item(91, 11).Value2 = 11850
$ws.Cells.Item(91, 12).Value2 = 0
$ws.Cells.Item(91, 14).Value2 = -10368
$ws.Cells.Item(95, 8).Value2 = 18400
$ws.Cells.Item(95, 10).Value2 = 0
$ws.Cells.Item(95, 12).Value2 = 0
$ws.Cells.Item(114, 8).Value2 = 429
$ws.Cells.Item(114, 9).Value2 = 515.5
$ws.Cells.Item(114, 10).Value2 = 394.4
$ws.Cells.Item(114, 11).Value2 = 1546.5
$ws.Cells.Item(114, 12).Value2 = 1183.2
$ws.Cells.Item(114, 13).Value2 = 1707.5
$ws.Cells.Item(114, 14).Value2 = -7691.2
$ws.Cells.Item(128, 8).Value2 = 436628.66
$ws.Cells.Item(128, 9).Value2 = 436628.66
$ws.Cells.Item(128, 11).Value2 = 1309885.98
$ws.Cells.Item(128, 13).Value2 = -1304905.98
$ws.Cells.Item(137, 8).Value2 = 2974.5454
$ws.Cells.Item(137, 9).Value2 = 780
$ws.Cells.Item(137, 11).Value2 = 2340
$ws.Cells.Item(137, 13).Value2 = 2760
$ws.Cells.Item(140, 8).Value2 = 1383.5
$ws.Cells.Item(140, 9).Value2 = 1383.5
$ws.Cells.Item(140, 11).Value2 = 4150.5
$ws.Cells.Item(140, 13).Value2 = 1029.5
$ws.Cells.Item(88, 13).ClearContents()
$ws.Cells.Item(91, 13).ClearContents()
$ws.Cells.Item(95, 14).ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(9, 8).Value2 = 959.75
$ws.Cells.Item(9, 10).Value2 = 696.1667
$ws.Cells.Item(9, 12).Value2 = 696.1667
$ws.Cells.Item(9, 14).Value2 = -1036.1667
$ws.Cells.Item(49, 8).Value2 = 59895
$ws.Cells.Item(49, 10).Value2 = 59895
$ws.Cells.Item(49, 12).Value2 = 59895
$ws.Cells.Item(49, 14).Value2 = -60263
$ws.Cells.Item(113, 8).Value2 = 23237.125
$ws.Cells.Item(113, 9).Value2 = 1756.6666
$ws.Cells.Item(113, 10).Value2 = 50854.855
$ws.Cells.Item(113, 11).Value2 = 1756.6666
$ws.Cells.Item(113, 12).Value2 = 50854.855
$ws.Cells.Item(113, 13).Value2 = 413.3334
$ws.Cells.Item(113, 14).Value2 = -55194.855
$ws.Cells.Item(122, 8).Value2 = 6505.125
$ws.Cells.Item(122, 9).Value2 = 7412.1304
$ws.Cells.Item(122, 10).Value2 = 4187.222
$ws.Cells.Item(122, 11).Value2 = 22236.3912
$ws.Cells.Item(122, 12).Value2 = 12561.666
$ws.Cells.Item(122, 13).Value2 = -19786.3912
$ws.Cells.Item(122, 14).Value2 = -17461.666

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value2 = 11371.131
$ws.Cells.Item(40, 9).Value2 = 10560.45
$ws.Cells.Item(40, 11).Value2 = 10560.45
$ws.Cells.Item(40, 13).Value2 = -10424.45
$ws.Cells.Item(122, 8).Value2 = 5995.9165
$ws.Cells.Item(122, 9).Value2 = 3993.625
$ws.Cells.Item(122, 10).Value2 = 10000.5
$ws.Cells.Item(122, 11).Value2 = 11980.875
$ws.Cells.Item(122, 12).Value2 = 30001.5
$ws.Cells.Item(122, 13).Value2 = -9530.875
$ws.Cells.Item(122, 14).Value2 = -34901.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(52, 8).Value2 = 12582.1
$ws.Cells.Item(52, 9).Value2 = 4503.6665
$ws.Cells.Item(52, 10).Value2 = 24699.75
$ws.Cells.Item(52, 11).Value2 = 4503.6665
$ws.Cells.Item(52, 12).Value2 = 24699.75
$ws.Cells.Item(52, 13).Value2 = -4277.6665
$ws.Cells.Item(52, 14).Value2 = -25151.75
$ws.Cells.Item(62, 8).Value2 = 31400.6
$ws.Cells.Item(62, 9).Value2 = 66500
$ws.Cells.Item(62, 10).Value2 = 8001
$ws.Cells.Item(62, 11).Value2 = 66500
$ws.Cells.Item(62, 12).Value2 = 8001
$ws.Cells.Item(62, 13).Value2 = -65876
$ws.Cells.Item(62, 14).Value2 = -9249
$ws.Cells.Item(65, 8).Value2 = 31400.6
$ws.Cells.Item(65, 9).Value2 = 66500
$ws.Cells.Item(65, 10).Value2 = 8001
$ws.Cells.Item(65, 11).Value2 = 332500
$ws.Cells.Item(65, 12).Value2 = 40005
$ws.Cells.Item(65, 13).Value2 = -329380
$ws.Cells.Item(65, 14).Value2 = -46245
$ws.Cells.Item(100, 8).Value2 = 969.05
$ws.Cells.Item(100, 9).Value2 = 881.125
$ws.Cells.Item(100, 11).Value2 = 1762.25
$ws.Cells.Item(100, 13).Value2 = -1221.25
